$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.956.09'
$ws.Range("E2").Value = '  +0.12%  '

$ws.Range("D3").Value = '2.668.82'
$ws.Range("E3").Value = '  +2.30%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.55'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.52%  '

$ws.Range("E7").Value = '  +0.20%  '

$ws.Range("E8").Value = '  -0.31%  '

$ws.Range("E9").Value = '  +0.71%  '

$ws.Range("E10").Value = '  +0.45%  '

$ws.Range("E11").Value = '  +2.83%  '

$ws.Range("E12").Value = '  -0.07%  '

$ws.Range("D13").Value = '3.120.27'
$ws.Range("E13").Value = '  +1.60%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.18'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +11.63%  '

$ws.Range("D15").Value = '60.929.53'
$ws.Range("E15").Value = '  +0.18%  '

$ws.Range("E16").Value = '  +0.55%  '

$ws.Range("D17").Value = '2.660.08'
$ws.Range("E17").Value = '  +1.45%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.59'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.28%  '

$ws.Range("E19").Value = '  +1.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '350.83'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.14%  '

$ws.Range("E21").Value = '  -0.45%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.529'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.69%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.99'
$ws.Range("D24").Style = "Normal"

$ws.Range("E25").Value = '  +0.77%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.996'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.49%  '

$ws.Range("E28").Value = '  +9.37%  '

$ws.Range("D29").Value = '0.0₃0811'
$ws.Range("E29").Value = '  +0.84%  '

$ws.Range("E30").Value = '  +7.52%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '163.76'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.77%  '

$ws.Range("E33").Value = '  +1.75%  '

$ws.Range("E34").Value = '  +6.69%  '

$ws.Range("E35").Value = '  +1.78%  '

$ws.Range("E36").Value = '  +7.03%  '

$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '339.48'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.97%  '

$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.66'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.50%  '

$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.09'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.10%  '

$ws.Range("B40").Value = 'SuiNetwork'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.913'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.11%  '

$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.66%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.27'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.74%  '

$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.624'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.63%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.36'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.72%  '

$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0565'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.46%  '

$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0250'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.75%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '133.16'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.86%  '

$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0997'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.01%  '

$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.62'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.47%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.46%  '

$ws.Range("D51").Value = '2.097.46'
$ws.Range("E51").Value = '  +3.41%  '
